# Update "想去人数" (interest count) figures on both the "展览" and
# "全部类型" worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3443
$ws1.Range("F3").Value = 26
$ws1.Range("F5").Value = 1741
$ws1.Range("F6").Value = 94
$ws1.Range("F7").Value = 339

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3443
$ws4.Range("F3").Value = 26
$ws4.Range("F5").Value = 1742
$ws4.Range("F6").Value = 94
$ws4.Range("F8").Value = 339
